# Realestate Update resale numbers 2024-01-02 09:33
# Append a new data row (row 5) to the CityResaleNum sheet with the
# 2024-01-02 09:33:49 resale snapshot values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 5

# Columns A (date-like "2024-01-02") and D ("00") would otherwise be
# auto-converted by Excel into a date serial / plain number, so force
# them to text first, matching how the other rows store these as text.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2024-01-02"

$ws.Cells.Item($row, 2).Value = "09:33:49"
$ws.Cells.Item($row, 3).Value = "Tuesday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "00"

$ws.Cells.Item($row, 5).Value = 140083
$ws.Cells.Item($row, 6).Value = 142925
$ws.Cells.Item($row, 7).Value = 171212
$ws.Cells.Item($row, 8).Value = 145338
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 116484
$ws.Cells.Item($row, 11).Value = 223673
$ws.Cells.Item($row, 12).Value = 247445
$ws.Cells.Item($row, 13).Value = 183416
$ws.Cells.Item($row, 14).Value = 109760
$ws.Cells.Item($row, 15).Value = 39603
$ws.Cells.Item($row, 16).Value = 30586
$ws.Cells.Item($row, 17).Value = 71574
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 41267
$ws.Cells.Item($row, 20).Value = -1
